# correction S/LFM mapping Italy commercial
#
# The sheet has one long multi-line string per column (B2=Offices, C2=Trade,
# D2=Hotels). Each line is "<share>% <taxonomy string>/<occupancy>".
# The "S/LFM+CDL/HBET:3-5/<occupancy>" line had been entered twice (with two
# different percentages) in each column; this fixes the mapping by removing
# the duplicate/incorrect entry and keeping the correct percentage on the
# remaining occurrence (for Offices, that means the surviving line's value
# becomes 0.6% instead of 0.0%).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-SLfmDuplicate {
    param(
        [string]$CellAddress,
        [string]$LineToDrop,
        [string]$SurvivingOldLine,
        [string]$SurvivingNewLine
    )

    $cell = $ws.Range($CellAddress)
    $text = $cell.Value2
    $lines = $text -split "`n"

    $dropIndex = -1
    for ($i = 0; $i -lt $lines.Length; $i++) {
        if ($lines[$i] -eq $LineToDrop) {
            $dropIndex = $i
            break
        }
    }
    if ($dropIndex -ge 0) {
        $newLines = @()
        for ($i = 0; $i -lt $lines.Length; $i++) {
            if ($i -ne $dropIndex) {
                $newLines += $lines[$i]
            }
        }
        $lines = $newLines
    }

    if ($SurvivingOldLine -ne $SurvivingNewLine) {
        for ($i = 0; $i -lt $lines.Length; $i++) {
            if ($lines[$i] -eq $SurvivingOldLine) {
                $lines[$i] = $SurvivingNewLine
            }
        }
    }

    $cell.Value = [string]::Join("`n", $lines)
}

# Offices (B2): drop the first "0.6% S/LFM+CDL/HBET:3-5/Offices" occurrence
# (it was mapped in the wrong spot) and correct the trailing occurrence from
# 0.0% to 0.6%.
Fix-SLfmDuplicate `
    "B2" `
    "0.6% S/LFM+CDL/HBET:3-5/Offices" `
    "0.0% S/LFM+CDL/HBET:3-5/Offices" `
    "0.6% S/LFM+CDL/HBET:3-5/Offices"

# Trade (C2): the 0.0% S/LFM line was duplicated; drop the first occurrence.
Fix-SLfmDuplicate `
    "C2" `
    "0.0% S/LFM+CDL/HBET:3-5/Trade" `
    "0.0% S/LFM+CDL/HBET:3-5/Trade" `
    "0.0% S/LFM+CDL/HBET:3-5/Trade"

# Hotels (D2): the " 0.0% S/LFM+CDL/HBET:3-5/Hotels" line was duplicated;
# drop the first occurrence (the correct one, " 2.0% ...", remains last).
Fix-SLfmDuplicate `
    "D2" `
    " 0.0% S/LFM+CDL/HBET:3-5/Hotels" `
    " 2.0% S/LFM+CDL/HBET:3-5/Hotels" `
    " 2.0% S/LFM+CDL/HBET:3-5/Hotels"

# Cosmetic view changes captured in the diff.
[void]$ws.Range("E2").Select()
$ws.Rows.Item(2).RowHeight = 350
